$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3.75
$ws.Range("U2").Value = 1.73
$ws.Range("V2").Value = 2
$ws.Range("AH2").Value = 17
$ws.Range("G3").Value = 1.73
$ws.Range("H3").Value = 3.8
$ws.Range("I3").Value = 4.5
$ws.Range("G4").Value = 2.63
$ws.Range("AD4").Value = 6.5
$ws.Range("AF4").Value = 41
$ws.Range("AJ4").Value = 10
$ws.Range("G5").Value = 1.75
$ws.Range("H5").Value = 3.4
$ws.Range("G6").Value = 2.55
$ws.Range("I6").Value = 2.45
$ws.Range("J6").Value = 3.2
$ws.Range("L6").Value = 3.1
$ws.Range("Q6").Value = 1.75
$ws.Range("R6").Value = 2.05
$ws.Range("S6").Value = 1.33
$ws.Range("T6").Value = 3.25
$ws.Range("X6").Value = 15
$ws.Range("AE6").Value = 12
$ws.Range("AO6").Value = 15
$ws.Range("AT6").Value = 3.25
$ws.Range("AX6").Value = 13
$ws.Range("J7").Value = 2.92
$ws.Range("K7").Value = 1.98
$ws.Range("L7").Value = 3.65
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 7.1
$ws.Range("Q7").Value = 2.2
$ws.Range("S7").Value = 1.44
$ws.Range("T7").Value = 2.4
$ws.Range("Z7").Value = 25
$ws.Range("AB7").Value = 35
$ws.Range("AC7").Value = 7.2
$ws.Range("AD7").Value = 5.8
$ws.Range("AI7").Value = 14.5
$ws.Range("AN7").Value = 4.1
$ws.Range("AO7").Value = 12.5
$ws.Range("AP7").Value = 21
$ws.Range("AQ7").Value = 50
$ws.Range("AR7").Value = 90
$ws.Range("AS7").Value = 300
$ws.Range("AT7").Value = 2.37
$ws.Range("AU7").Value = 7.2
$ws.Range("AV7").Value = 70
$ws.Range("AW7").Value = 4.8
$ws.Range("AX7").Value = 17
$ws.Range("AY7").Value = 26
$ws.Range("AZ7").Value = 90
$ws.Range("BA7").Value = 150
$ws.Range("BB7").Value = 400
$ws.Range("G9").Value = 1.33
$ws.Range("H9").Value = 4.6
$ws.Range("I9").Value = 7.2
$ws.Range("J9").Value = 1.78
$ws.Range("L9").Value = 6.3
$ws.Range("M9").Value = 1.02
$ws.Range("N9").Value = 13.3
$ws.Range("O9").Value = 1.14
$ws.Range("P9").Value = 4.4
$ws.Range("Q9").Value = 1.53
$ws.Range("R9").Value = 2.2
$ws.Range("S9").Value = 1.26
$ws.Range("T9").Value = 3.5
$ws.Range("U9").Value = 1.83
$ws.Range("V9").Value = 1.93
$ws.Range("W9").Value = 6.8
$ws.Range("X9").Value = 6
$ws.Range("Y9").Value = 7.2
$ws.Range("Z9").Value = 7.4
$ws.Range("AB9").Value = 19
$ws.Range("AC9").Value = 14.5
$ws.Range("AD9").Value = 8.25
$ws.Range("AE9").Value = 15.5
$ws.Range("AF9").Value = 60
$ws.Range("AG9").Value = 350
$ws.Range("AH9").Value = 18
$ws.Range("AI9").Value = 40
$ws.Range("AJ9").Value = 18.5
$ws.Range("AK9").Value = 120
$ws.Range("AL9").Value = 60
$ws.Range("AM9").Value = 50
$ws.Range("AN9").Value = 3.25
$ws.Range("AO9").Value = 5.9
$ws.Range("AP9").Value = 15
$ws.Range("AQ9").Value = 15.5
$ws.Range("AR9").Value = 40
$ws.Range("AS9").Value = 175
$ws.Range("AT9").Value = 3.25
$ws.Range("AU9").Value = 8
$ws.Range("AV9").Value = 70
$ws.Range("AW9").Value = 8.5
$ws.Range("AX9").Value = 40
$ws.Range("AY9").Value = 37
$ws.Range("AZ9").Value = 250
$ws.Range("BA9").Value = 250
$ws.Range("BB9").Value = 450
$ws.Range("H10").Value = 4.4
$ws.Range("I10").Value = 1.4
$ws.Range("J10").Value = 5.7
$ws.Range("K10").Value = 2.4
$ws.Range("L10").Value = 1.87
$ws.Range("M10").Value = 1.02
$ws.Range("N10").Value = 13.4
$ws.Range("O10").Value = 1.2
$ws.Range("P10").Value = 3.68
$ws.Range("Q10").Value = 1.55
$ws.Range("R10").Value = 2.15
$ws.Range("S10").Value = 1.33
$ws.Range("T10").Value = 3.14
$ws.Range("U10").Value = 1.81
$ws.Range("V10").Value = 1.95
$ws.Range("W10").Value = 15
$ws.Range("AB10").Value = 45
$ws.Range("AC10").Value = 14
$ws.Range("AD10").Value = 7.8
$ws.Range("AE10").Value = 14.5
$ws.Range("AF10").Value = 55
$ws.Range("AG10").Value = 350
$ws.Range("AH10").Value = 6.8
$ws.Range("AI10").Value = 6.2
$ws.Range("AJ10").Value = 7.1
$ws.Range("AK10").Value = 8
$ws.Range("AL10").Value = 9.25
$ws.Range("AM10").Value = 19
$ws.Range("AN10").Value = 7.6
$ws.Range("AO10").Value = 35
$ws.Range("AP10").Value = 35
$ws.Range("AR10").Value = 250
$ws.Range("AS10").Value = 450
$ws.Range("AT10").Value = 3.2
$ws.Range("AU10").Value = 7.9
$ws.Range("AV10").Value = 70
$ws.Range("AW10").Value = 3.3
$ws.Range("AX10").Value = 6.3
$ws.Range("AY10").Value = 15.5
$ws.Range("AZ10").Value = 17.5
$ws.Range("BA10").Value = 40
$ws.Range("BB10").Value = 200
$ws.Range("P11").Value = 8.8
$ws.Range("Q11").Value = 1.22
$ws.Range("R11").Value = 3.5
$ws.Range("U11").Value = 2.84
$ws.Range("V11").Value = 1.4
$ws.Range("H12").Value = 3.7
$ws.Range("K12").Value = 2.1
$ws.Range("O12").Value = 1.3
$ws.Range("P12").Value = 3.4
$ws.Range("R12").Value = 1.8
$ws.Range("U12").Value = 2
$ws.Range("V12").Value = 1.73
$ws.Range("W12").Value = 13
$ws.Range("AB12").Value = 41
$ws.Range("AC12").Value = 9
$ws.Range("AE12").Value = 19
$ws.Range("AF12").Value = 51
$ws.Range("AH12").Value = 6.5
$ws.Range("AI12").Value = 7.5
$ws.Range("AM12").Value = 29
$ws.Range("AR12").Value = 126
$ws.Range("AU12").Value = 9
$ws.Range("AW12").Value = 3.6
$ws.Range("AX12").Value = 8.5
$ws.Range("AY12").Value = 21
$ws.Range("S13").Value = 1.4
$ws.Range("T13").Value = 2.75
$ws.Range("U13").Value = 1.73
$ws.Range("V13").Value = 2
$ws.Range("W13").Value = 8.5
$ws.Range("AC13").Value = 10
$ws.Range("AG13").Value = 201
$ws.Range("AO13").Value = 13
$ws.Range("AT13").Value = 2.75
$ws.Range("BA13").Value = 67
$ws.Range("BB13").Value = 151
$ws.Range("Q14").Value = 1.95
$ws.Range("R14").Value = 1.9
$ws.Range("H15").Value = 3.2
$ws.Range("I15").Value = 3.6
$ws.Range("M15").Value = 1.08
$ws.Range("N15").Value = 8
$ws.Range("S15").Value = 1.5
$ws.Range("T15").Value = 2.5
$ws.Range("Y15").Value = 9.5
$ws.Range("AB15").Value = 34
$ws.Range("AC15").Value = 8
$ws.Range("AG15").Value = 351
$ws.Range("AT15").Value = 2.5
$ws.Range("AV15").Value = 67
$ws.Range("G16").Value = 2.55
$ws.Range("I16").Value = 2.88
$ws.Range("J16").Value = 3.4
$ws.Range("X16").Value = 11
$ws.Range("Y16").Value = 11
$ws.Range("Z16").Value = 26
$ws.Range("AH16").Value = 7
$ws.Range("AJ16").Value = 11
$ws.Range("AK16").Value = 29
$ws.Range("AL16").Value = 26
$ws.Range("AN16").Value = 4.5
$ws.Range("AW16").Value = 4.75
$ws.Range("AX16").Value = 17
$ws.Range("AY16").Value = 29
$ws.Range("AZ16").Value = 51
$ws.Range("K17").Value = 2.55
$ws.Range("L17").Value = 3.7
$ws.Range("O17").Value = 1.11
$ws.Range("P17").Value = 5.6
$ws.Range("Q17").Value = 1.35
$ws.Range("R17").Value = 2.95
$ws.Range("T17").Value = 3.9
$ws.Range("U17").Value = 1.37
$ws.Range("V17").Value = 2.87
$ws.Range("AA17").Value = 12
$ws.Range("AC17").Value = 10.25
$ws.Range("AD17").Value = 9.5
$ws.Range("AF17").Value = 28
$ws.Range("AG17").Value = 120
$ws.Range("AH17").Value = 20
$ws.Range("AI17").Value = 27
$ws.Range("AM17").Value = 23
$ws.Range("AN17").Value = 4.45
$ws.Range("AQ17").Value = 24
$ws.Range("AT17").Value = 3.9
$ws.Range("AU17").Value = 6.2
$ws.Range("AV17").Value = 32
$ws.Range("AW17").Value = 6.2
$ws.Range("AX17").Value = 17.5
$ws.Range("AY17").Value = 17.5
$ws.Range("AZ17").Value = 70
$ws.Range("BA17").Value = 70
$ws.Range("BB17").Value = 150
